# Update "想去人数" (F column) figures across the sheets to reflect the
# newly generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Changes = @{
    4  = 5893
    5  = 5893
    7  = 2942
    8  = 1267
    10 = 434
    13 = 690
    14 = 191
    15 = 4244
    16 = 4244
    18 = 85
    19 = 101
    22 = 60
    23 = 6381
    24 = 6381
    29 = 217
    30 = 443
    31 = 5924
    34 = 1847
    35 = 5910
    36 = 102
    40 = 287
    41 = 4012
    43 = 76
    45 = 2395
    50 = 293
    51 = 2034
}
foreach ($row in $sheet1Changes.Keys) {
    $ws1.Range("F$row").Value = $sheet1Changes[$row]
}

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Changes = @{
    3  = 186
    5  = 93
    11 = 15
}
foreach ($row in $sheet2Changes.Keys) {
    $ws2.Range("F$row").Value = $sheet2Changes[$row]
}

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1408

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Changes = @{
    2  = 1408
    4  = 5893
    5  = 5893
    7  = 2942
    8  = 1267
    9  = 434
    12 = 186
    13 = 191
    14 = 4244
    15 = 4244
    17 = 85
    18 = 101
    21 = 60
    22 = 6381
    23 = 6381
    27 = 217
    28 = 93
    29 = 5924
    33 = 1847
    35 = 5910
    36 = 102
    40 = 4012
    42 = 76
    46 = 2395
    51 = 293
    52 = 15
}
foreach ($row in $sheet4Changes.Keys) {
    $ws4.Range("F$row").Value = $sheet4Changes[$row]
}
